$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# mis-interpreted by Excel as numbers (losing exact text representation,
# e.g. trailing zeros or precision), so they remain exact strings like
# the original inline-string cells.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.857.47"
$ws.Range("E2").Value = "  +4.48%  "
$ws.Range("D3").Value = "3.514.74"
$ws.Range("E3").Value = "  +2.38%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "594.78"
$ws.Range("E5").Value = "  +3.90%  "
$ws.Range("D6").Value = "169.47"
$ws.Range("E6").Value = "  +5.98%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "3.514.80"
$ws.Range("E8").Value = "  +2.33%  "
$ws.Range("D9").Value = "0.572"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "7.29"
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("D11").Value = "0.126"
$ws.Range("E11").Value = "  +5.13%  "
$ws.Range("D12").Value = "0.439"
$ws.Range("E12").Value = "  +3.05%  "
$ws.Range("D13").Value = "4.119.99"
$ws.Range("E13").Value = "  +2.39%  "
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "28.20"
$ws.Range("E15").Value = "  +3.42%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000179"
$ws.Range("E16").Value = "  +2.98%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "66.817.73"
$ws.Range("E17").Value = "  +4.28%  "
$ws.Range("D18").Value = "3.488.68"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("D19").Value = "6.34"
$ws.Range("E19").Value = "  +3.84%  "
$ws.Range("D20").Value = "14.07"
$ws.Range("E20").Value = "  +3.03%  "
$ws.Range("D21").Value = "391.39"
$ws.Range("E21").Value = "  +1.96%  "
$ws.Range("D22").Value = "7.99"
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("D23").Value = "73.58"
$ws.Range("E23").Value = "  +2.91%  "
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").Value = "0.0000127"
$ws.Range("E25").Value = "  +8.33%  "
$ws.Range("D26").Value = "0.534"
$ws.Range("E26").Value = "  +3.21%  "
$ws.Range("D27").Value = "10.15"
$ws.Range("E27").Value = "  +3.92%  "
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "6.43"
$ws.Range("E30").Value = "  +5.90%  "
$ws.Range("E31").Value = "  +4.64%  "
$ws.Range("E32").Value = "  +3.01%  "
$ws.Range("D33").Value = "23.62"
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("D34").Value = "7.46"
$ws.Range("E34").Value = "  +6.69%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "1.62"
$ws.Range("E36").Value = "  +6.34%  "
$ws.Range("D37").Value = "161.36"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").Value = "0.895"
$ws.Range("E38").Value = "  +5.29%  "
$ws.Range("E39").Value = "  +3.32%  "
$ws.Range("D40").Value = "0.0748"
$ws.Range("E40").Value = "  +3.10%  "
$ws.Range("D41").Value = "26.72"
$ws.Range("E41").Value = "  +2.25%  "
$ws.Range("D42").Value = "4.64"
$ws.Range("E42").Value = "  +5.80%  "
$ws.Range("D43").Value = "2.849.78"
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").Value = "6.67"
$ws.Range("E44").Value = "  +3.61%  "
$ws.Range("D45").Value = "43.51"
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("D46").Value = "26.47"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D48").Value = "2.54"
$ws.Range("E48").Value = "  +3.04%  "
$ws.Range("D49").Value = "354.40"
$ws.Range("E49").Value = "  +5.40%  "
$ws.Range("E50").Value = "  +2.91%  "
$ws.Range("D51").Value = "33.67"
$ws.Range("E51").Value = "  +12.81%  "
